# Update cryptos.xlsx price/volume data to the latest snapshot.
# Column D (Price) values that look like plain numbers are written with a
# leading apostrophe so Excel stores/keeps them as text (matching the
# original inlineStr cell type) instead of silently converting them to
# numeric cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.991.26'
$ws.Range("E2").Value = '  +0.41%  '
$ws.Range("D3").Value = '3.319.99'
$ws.Range("E3").Value = '  +1.49%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '''562.05'
$ws.Range("E5").Value = '  +1.02%  '
$ws.Range("D6").Value = '''185.80'
$ws.Range("E6").Value = '  +0.83%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.312.93'
$ws.Range("E8").Value = '  +1.51%  '
$ws.Range("D9").Value = '''0.574'
$ws.Range("E9").Value = '  -2.45%  '
$ws.Range("D10").Value = '''0.176'
$ws.Range("E10").Value = '  -4.89%  '
$ws.Range("E11").Value = '  -1.00%  '
$ws.Range("D12").Value = '''46.02'
$ws.Range("E12").Value = '  -2.30%  '
$ws.Range("D13").Value = '''0.0000264'
$ws.Range("E13").Value = '  -0.41%  '
$ws.Range("D14").Value = '3.851.55'
$ws.Range("E14").Value = '  +1.55%  '
$ws.Range("E15").Value = '  -2.14%  '
$ws.Range("D16").Value = '''581.81'
$ws.Range("E16").Value = '  -8.55%  '
$ws.Range("D17").Value = '66.029.73'
$ws.Range("E17").Value = '  +0.53%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.330.84'
$ws.Range("E18").Value = '  +2.14%  '
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").Value = '''0.117'
$ws.Range("E19").Value = '  +0.32%  '
$ws.Range("D20").Value = '''17.69'
$ws.Range("E20").Value = '  -1.81%  '
$ws.Range("D21").Value = '''10.89'
$ws.Range("E21").Value = '  -3.88%  '
$ws.Range("D22").Value = '''0.895'
$ws.Range("E22").Value = '  -0.67%  '
$ws.Range("D23").Value = '''17.69'
$ws.Range("E23").Value = '  -2.88%  '
$ws.Range("D24").Value = '''5.04'
$ws.Range("E24").Value = '  +2.63%  '
$ws.Range("D25").Value = '''97.74'
$ws.Range("E25").Value = '  -9.22%  '
$ws.Range("D26").Value = '''3.97'
$ws.Range("E26").Value = '  +0.28%  '
$ws.Range("B27").Value = 'ImmutableX'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D27").Value = '''2.71'
$ws.Range("E27").Value = '  +1.42%  '
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").Value = '''9.38'
$ws.Range("E28").Value = '  -1.84%  '
$ws.Range("B29").Value = 'Filecoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D29").Value = '''8.46'
$ws.Range("E29").Value = '  -2.12%  '
$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").Value = '''30.67'
$ws.Range("E30").Value = '  +1.67%  '
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").Value = '''6.69'
$ws.Range("E31").Value = '  +7.16%  '
$ws.Range("B32").Value = 'dogwifhat'
$ws.Range("C32").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D32").Value = '''3.69'
$ws.Range("E32").Value = '  -6.67%  '
$ws.Range("B33").Value = 'Bittensor'
$ws.Range("C33").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D33").Value = '''563.65'
$ws.Range("E33").Value = '  +9.47%  '
$ws.Range("B34").Value = 'Cosmos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D34").Value = '''10.82'
$ws.Range("E34").Value = '  -1.77%  '
$ws.Range("B35").Value = 'Maker'
$ws.Range("C35").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D35").Value = '3.774.10'
$ws.Range("E35").Value = '  +0.58%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = '''0.103'
$ws.Range("E36").Value = '  -1.45%  '
$ws.Range("B37").Value = 'Dai'
$ws.Range("C37").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D37").Value = '''1.00'
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("B38").Value = 'OKB'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D38").Value = '''55.69'
$ws.Range("E38").Value = '  -3.24%  '
$ws.Range("B39").Value = 'InjectiveProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D39").Value = '''33.28'
$ws.Range("E39").Value = '  +0.98%  '
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").Value = '''0.127'
$ws.Range("E40").Value = '  -2.71%  '
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").Value = '''3.15'
$ws.Range("E41").Value = '  -6.94%  '
$ws.Range("B42").Value = 'PEPE'
$ws.Range("C42").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D42").Value = '0.0₃0682'
$ws.Range("E42").Value = '  -7.01%  '
$ws.Range("B43").Value = 'Fetch.AI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D43").Value = '''2.59'
$ws.Range("E43").Value = '  -4.83%  '
$ws.Range("B44").Value = 'ApeXProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D44").Value = '''3.34'
$ws.Range("E44").Value = '  +4.30%  '
$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").Value = '''0.332'
$ws.Range("E45").Value = '  -0.94%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = '''0.0409'
$ws.Range("E46").Value = '  -0.69%  '
$ws.Range("D47").Value = '''3.06'
$ws.Range("E47").Value = '  -11.02%  '
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").Value = '''0.126'
$ws.Range("E48").Value = '  -2.01%  '
$ws.Range("B49").Value = 'FirstDigitalUSD'
$ws.Range("C49").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D49").Value = '''1.00'
$ws.Range("E49").Value = '  +0.11%  '
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").Value = '''2.51'
$ws.Range("E50").Value = '  -3.53%  '
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").Value = '''127.07'
$ws.Range("E51").Value = '  +4.38%  '
